# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.946.86"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +1.41%  "

# Row 3
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.663.45"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -0.92%  "

# Row 4
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.93%  "

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.71"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +3.38%  "

# Row 6
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.28%  "

# Row 7
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3642"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -1.22%  "

# Row 8
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.23"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -1.73%  "

# Row 9
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3286"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -2.08%  "

# Row 10
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.144"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -2.91%  "

# Row 11
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07097"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -2.68%  "

# Row 12
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -0.51%  "

# Row 13
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.088"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -1.30%  "

# Row 14
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.74"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -3.36%  "

# Row 15
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.665.40"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -1.24%  "

# Row 16
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.644"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -2.52%  "

# Row 17
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001054"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -3.76%  "

# Row 18
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06667"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +0.59%  "

# Row 19
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9976"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -0.39%  "

# Row 20
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "79.67"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -2.02%  "

# Row 21
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.944"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -4.02%  "

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.83"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -5.65%  "

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.64"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.66%  "

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.889.10"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +1.00%  "

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.438"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -0.06%  "

# Row 26
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.421"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -9.70%  "

# Row 27
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.07"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +0.99%  "

# Row 28
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.71"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -5.48%  "

# Row 29
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.237"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +3.88%  "

# Row 30
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.849.99"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -1.15%  "

# Row 31
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.18"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -2.89%  "

# Row 32
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.115"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -1.66%  "

# Row 33
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.879"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -9.30%  "

# Row 34
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08520"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -0.75%  "

# Row 35
$ws.Range("E35").Value = "  -3.31%  "

# Row 36
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.32"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -6.68%  "

# Row 37
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.285"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +3.48%  "

# Row 38
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.246"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -2.87%  "

# Row 39
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02271"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -2.38%  "

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06103"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -5.25%  "

# Row 41
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.351"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -5.35%  "

# Row 42
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2083"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -3.24%  "

# Row 43
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9982"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -0.36%  "

# Row 44
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5981"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -4.06%  "

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.827"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +1.46%  "

# Row 46
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.82"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -3.90%  "

# Row 47
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5673"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -4.06%  "

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.46"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +0.59%  "

# Row 49
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.970"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -3.94%  "

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07035"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -1.79%  "

# Row 51
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.203"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +1.16%  "
